$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the Price column as text so values like "5.80" or "60.836.01" are not
# reinterpreted/rounded as numbers by Excel's auto-detection.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.836.01"
$ws.Range("E2").Value = "  -1.39%  "

$ws.Range("D3").Value = "3.372.51"
$ws.Range("E3").Value = "  -0.59%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "572.49"
$ws.Range("E5").Value = "  -0.57%  "

$ws.Range("D6").Value = "137.03"
$ws.Range("E6").Value = "  +0.65%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "3.370.47"

$ws.Range("E9").Value = "  -0.85%  "

$ws.Range("D10").Value = "7.66"
$ws.Range("E10").Value = "  +2.90%  "

$ws.Range("E11").Value = "  -2.04%  "

$ws.Range("D12").Value = "0.385"
$ws.Range("E12").Value = "  -1.58%  "

$ws.Range("D13").Value = "3.950.65"
$ws.Range("E13").Value = "  -0.60%  "

$ws.Range("E14").Value = "  +0.50%  "

$ws.Range("D15").Value = "25.96"
$ws.Range("E15").Value = "  +2.47%  "

$ws.Range("E16").Value = "  -2.26%  "

$ws.Range("D17").Value = "3.376.48"
$ws.Range("E17").Value = "  -0.71%  "

$ws.Range("D18").Value = "61.029.33"
$ws.Range("E18").Value = "  -1.20%  "

$ws.Range("D19").Value = "13.95"
$ws.Range("E19").Value = "  -1.68%  "

$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "5.80"
$ws.Range("E20").Value = "  -0.72%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "9.42"
$ws.Range("E21").Value = "  -0.23%  "

$ws.Range("D22").Value = "373.99"
$ws.Range("E22").Value = "  -3.34%  "

$ws.Range("E23").Value = "  -2.69%  "

$ws.Range("D24").Value = "3.519.43"
$ws.Range("E24").Value = "  -0.57%  "

$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.13%  "

$ws.Range("E26").Value = "  -0.91%  "

$ws.Range("D27").Value = "71.11"
$ws.Range("E27").Value = "  -0.12%  "

$ws.Range("D28").Value = "1.78"
$ws.Range("E28").Value = "  +12.31%  "

$ws.Range("E29").Value = "  +10.35%  "

$ws.Range("E30").Value = "  -2.12%  "

$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.37%  "

$ws.Range("D32").Value = "8.11"
$ws.Range("E32").Value = "  -1.73%  "

$ws.Range("E33").Value = "  -0.81%  "

$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("D35").Value = "23.66"
$ws.Range("E35").Value = "  +0.79%  "

$ws.Range("D36").Value = "5.18"
$ws.Range("E36").Value = "  -4.07%  "

$ws.Range("E37").Value = "  -1.28%  "

$ws.Range("E38").Value = "  -0.57%  "

$ws.Range("D39").Value = "164.84"
$ws.Range("E39").Value = "  +1.01%  "

$ws.Range("D40").Value = "0.0762"
$ws.Range("E40").Value = "  -3.02%  "

$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("D42").Value = "0.775"
$ws.Range("E42").Value = "  -0.91%  "

$ws.Range("D43").Value = "41.61"
$ws.Range("E43").Value = "  -0.21%  "

$ws.Range("E44").Value = "  -4.74%  "

$ws.Range("D45").Value = "4.39"
$ws.Range("E45").Value = "  -1.11%  "

$ws.Range("D46").Value = "1.19"
$ws.Range("E46").Value = "  -2.63%  "

$ws.Range("D47").Value = "24.25"
$ws.Range("E47").Value = "  -1.45%  "

$ws.Range("D48").Value = "2.457.81"
$ws.Range("E48").Value = "  +4.23%  "

$ws.Range("D49").Value = "6.79"
$ws.Range("E49").Value = "  -2.12%  "

$ws.Range("D50").Value = "22.95"
$ws.Range("E50").Value = "  -1.28%  "

$ws.Range("D51").Value = "2.41"
$ws.Range("E51").Value = "  +4.37%  "
